# Change 2 on excel: add a second row of data (A2 = "Hi ", B2 = "World")
# and move the active selection to B2, matching the target workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Hi "
$ws.Range("B2").Value = "World"

$ws.Range("B2").Select()
